$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 260
$firstRow = 2

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
